$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

$ws.Cells.Item($row, 1).Value = 8
$ws.Cells.Item($row, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 45191
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item($row - 1, 4).NumberFormat
$ws.Cells.Item($row, 5).Value = 4
$ws.Cells.Item($row, 6).Value = 100112026
$ws.Cells.Item($row, 7).Value = "Haba"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 400
$ws.Cells.Item($row, 11).Value = 9000
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 9500
$ws.Cells.Item($row, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia del Elquí"
$ws.Cells.Item($row, 16).Value = 380
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
